$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("B2").Value = 0.0000008912037037037037
$ws.Range("C2").Value = 0.0000290625

# Data for the new rows 3-11: A, B, C
$data = @(
    @(10,  0,                              0.000001805555555555556),
    @(50,  0.00000002314814814814815,      0.00005068287037037037),
    @(100, 0.00000008101851851851852,      0.0002020833333333333),
    @(10,  0,                              0.000002881944444444444),
    @(20,  0,                              0.00000974537037037037),
    @(10,  0,                              0.000005578703703703703),
    @(20,  0,                              0.000006562499999999999),
    @(5,   0,                              0.000001365740740740741),
    @(10,  0.0000000009027777777777777,    0.000002091493055555555)
)

$row = 3
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]

    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 2).NumberFormat = "[hh]:mm:ss"

    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 3).NumberFormat = "[hh]:mm:ss"

    $row++
}
